# Updated cryptos list on Fri Jan 26 09:57:24 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.232.07"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "2.217.49"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "297.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.08"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.473"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.75"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.16"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "2.563.13"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.91"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "2.209.86"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.738"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").Value = "40.207.70"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.80"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.82"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.33"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.37"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.46"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.31"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.99"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.99"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0716"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.70%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.114"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.74"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.63"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.84"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "2.067.07"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.44"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.89%  "
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.01"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("E47").Value = "  +7.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.88"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -12.10%  "
$ws.Range("D49").Value = "2.434.66"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("E51").Value = "  +1.11%  "
